$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 71

# Values (copy the data row below the header, matching the pattern used by prior rows)
$ws.Cells.Item($row, 1).Value = 70
$ws.Cells.Item($row, 2).Value = "algeria"
$ws.Cells.Item($row, 3).Value = "ligue-1"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45276.625
$ws.Cells.Item($row, 6).Value = "US Souf"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "ES Setif"
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 3.76
$ws.Cells.Item($row, 11).Value = "14/12/2023 09:13"
$ws.Cells.Item($row, 12).Value = 3.07
$ws.Cells.Item($row, 13).Value = "16/12/2023 14:36"
$ws.Cells.Item($row, 14).Value = 2.86
$ws.Cells.Item($row, 15).Value = "14/12/2023 09:13"
$ws.Cells.Item($row, 16).Value = 2.92
$ws.Cells.Item($row, 17).Value = "16/12/2023 13:05"
$ws.Cells.Item($row, 18).Value = 2.09
$ws.Cells.Item($row, 19).Value = "14/12/2023 09:13"
$ws.Cells.Item($row, 20).Value = 2.57
$ws.Cells.Item($row, 21).Value = "16/12/2023 14:36"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-es-setif/zgNe6gJB/"

# Formatting: columns A and E carry specific cell styles (bordered/bold index
# for A, date-time number format for E) in every existing data row. Copy the
# formats from the row above (row 70) so the new row matches exactly.
$ws.Range("A70").Copy()
$ws.Range("A71").PasteSpecial(-4122)

$ws.Range("E70").Copy()
$ws.Range("E71").PasteSpecial(-4122)

$excel.CutCopyMode = 0
